{"js": "// Update three-digit x one-digit multiplication problems throughout the document's table.\nconst body = context.document.body;\nconst replacements = [\n  [\"346\u00d78=\", \"889\u00d77=\"],\n  [\"262\u00d73=\", \"172\u00d79=\"],\n  [\"739\u00d78=\", \"723\u00d78=\"],\n  [\"576\u00d75=\", \"540\u00d78=\"],\n  [\"702\u00d79=\", \"376\u00d77=\"],\n  [\"817\u00d77=\", \"778\u00d73=\"],\n  [\"417\u00d77=\", \"166\u00d72=\"],\n  [\"584\u00d72=\", \"630\u00d72=\"],\n  [\"513\u00d75=\", \"831\u00d77=\"],\n  [\"580\u00d74=\", \"422\u00d79=\"],\n  [\"978\u00d77=\", \"963\u00d79=\"],\n  [\"742\u00d72=\", \"659\u00d77=\"],\n  [\"437\u00d75=\", \"495\u00d72=\"],\n  [\"961\u00d72=\", \"432\u00d76=\"],\n  [\"151\u00d77=\", \"279\u00d77=\"],\n  [\"605\u00d75=\", \"609\u00d72=\"],\n  [\"212\u00d74=\", \"113\u00d74=\"],\n  [\"333\u00d73=\", \"201\u00d72=\"],\n  [\"182\u00d74=\", \"396\u00d73=\"],\n  [\"215\u00d79=\", \"325\u00d79=\"],\n  [\"810\u00d79=\", \"434\u00d72=\"],\n  [\"908\u00d76=\", \"510\u00d75=\"],\n  [\"684\u00d72=\", \"227\u00d77=\"],\n  [\"458\u00d79=\", \"928\u00d77=\"],\n  [\"909\u00d73=\", \"792\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update three-digit x one-digit multiplication problems throughout the document's table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"346\u00d78=\", \"889\u00d77=\"),\n    @(\"262\u00d73=\", \"172\u00d79=\"),\n    @(\"739\u00d78=\", \"723\u00d78=\"),\n    @(\"576\u00d75=\", \"540\u00d78=\"),\n    @(\"702\u00d79=\", \"376\u00d77=\"),\n    @(\"817\u00d77=\", \"778\u00d73=\"),\n    @(\"417\u00d77=\", \"166\u00d72=\"),\n    @(\"584\u00d72=\", \"630\u00d72=\"),\n    @(\"513\u00d75=\", \"831\u00d77=\"),\n    @(\"580\u00d74=\", \"422\u00d79=\"),\n    @(\"978\u00d77=\", \"963\u00d79=\"),\n    @(\"742\u00d72=\", \"659\u00d77=\"),\n    @(\"437\u00d75=\", \"495\u00d72=\"),\n    @(\"961\u00d72=\", \"432\u00d76=\"),\n    @(\"151\u00d77=\", \"279\u00d77=\"),\n    @(\"605\u00d75=\", \"609\u00d72=\"),\n    @(\"212\u00d74=\", \"113\u00d74=\"),\n    @(\"333\u00d73=\", \"201\u00d72=\"),\n    @(\"182\u00d74=\", \"396\u00d73=\"),\n    @(\"215\u00d79=\", \"325\u00d79=\"),\n    @(\"810\u00d79=\", \"434\u00d72=\"),\n    @(\"908\u00d76=\", \"510\u00d75=\"),\n    @(\"684\u00d72=\", \"227\u00d77=\"),\n    @(\"458\u00d79=\", \"928\u00d77=\"),\n    @(\"909\u00d73=\", \"792\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
